# Add "2022-Q3" fund-holdings data:
#  1. Update the "总计" (summary) sheet: insert a new 2022-Q3 row at the top
#     of the data (pushing every later quarter down one row, with a brand
#     new last row for 2021-Q2).
#  2. Insert a new worksheet named "2022-Q3" right after "总计" (i.e. as the
#     new second sheet, before "2022-Q2") containing the per-fund detail
#     rows for that quarter.
#  3. Restore "2021-Q2" as the selected/active sheet (it was the active tab
#     before the edit, and none of this edit is meant to change that).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "总计" summary sheet — shift the quarterly rows down by one and fill
#    in the new 2022-Q3 summary figures.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 5
$summary.Range("D2").Value = 0.87

$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 25
$summary.Range("D3").Value = 2.87

$summary.Range("B4").Value = "2022-Q1"
$summary.Range("C4").Value = 9
$summary.Range("D4").Value = 1.5

$summary.Range("B5").Value = "2021-Q4"
$summary.Range("C5").Value = 5
$summary.Range("D5").Value = 0.32

$summary.Range("B6").Value = "2021-Q3"
$summary.Range("C6").Value = 1
$summary.Range("D6").Value = 0.03

# Brand new row 7 (2021-Q2, previously the last row in the table). Clone
# the index column's formatting (bold / centered / bordered, style of the
# other A2:A6 index cells) onto the new A7 cell before writing its value.
$summary.Range("A2").Copy()
$summary.Range("A7").PasteSpecial(-4122)
$summary.Range("A7").Value = 5
$summary.Range("B7").Value = "2021-Q2"
$summary.Range("C7").Value = 1
$summary.Range("D7").Value = 0.02

# ---------------------------------------------------------------------
# 2. New "2022-Q3" worksheet with the per-fund breakdown. "2021-Q4" has
#    the exact same shape (header row + 5 funds, columns A-H) so clone it
#    as a template — this carries over the header/index-column styling —
#    then overwrite every cell with the 2022-Q3 figures.
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$template.Copy($wb.Worksheets.Item("2022-Q2"))
$q3 = $wb.Worksheets.Item("2021-Q4 (2)")
$q3.Name = "2022-Q3"

# Row 2 — 东方红创新趋势混合
$q3.Range("B2").Value = "'010699"
$q3.Range("C2").Value = "东方红创新趋势混合"
$q3.Range("D2").Value = "'26.13"
$q3.Range("E2").Value = "'71.44"
$q3.Range("F2").Value = "'2.48"
$q3.Range("G2").Value = "'0.6480"
$q3.Range("H2").Value = 9

# Row 3 — 东方红启华三年持有期混合A
$q3.Range("B3").Value = "'910021"
$q3.Range("C3").Value = "东方红启华三年持有期混合A"
$q3.Range("D3").Value = "'4.33"
$q3.Range("E3").Value = "'74.71"
$q3.Range("F3").Value = "'3.81"
$q3.Range("G3").Value = "'0.1650"
$q3.Range("H3").Value = 3

# Row 4 — 东方红启华三年持有期混合B
$q3.Range("B4").Value = "'011313"
$q3.Range("C4").Value = "东方红启华三年持有期混合B"
$q3.Range("D4").Value = "'0.90"
$q3.Range("E4").Value = "'74.71"
$q3.Range("F4").Value = "'3.81"
$q3.Range("G4").Value = "'0.0343"
$q3.Range("H4").Value = 3

# Row 5 — 民生加银新战略灵活配置混合A
$q3.Range("B5").Value = "'001352"
$q3.Range("C5").Value = "民生加银新战略灵活配置混合A"
$q3.Range("D5").Value = "'0.77"
$q3.Range("E5").Value = "'46.20"
$q3.Range("F5").Value = "'3.05"
$q3.Range("G5").Value = "'0.0235"
$q3.Range("H5").Value = 8

# Row 6 — 民生加银新战略灵活配置混合C
$q3.Range("B6").Value = "'011391"
$q3.Range("C6").Value = "民生加银新战略灵活配置混合C"
$q3.Range("D6").Value = "'0.02"
$q3.Range("E6").Value = "'46.20"
$q3.Range("F6").Value = "'3.05"
$q3.Range("G6").Value = "'0.0006"
$q3.Range("H6").Value = 8

# The leading "'" above forces the numeric-looking strings (fund codes,
# percentages, AUM figures) to stay text, same as every other quarter's
# sheet — but it also stamps a stray quote-prefix format on those cells.
# Strip that back off by re-pasting the plain (unstyled) format from the
# fund-name column, which was never touched by the quote-prefix trick.
$q3.Range("C2").Copy()
$q3.Range("B2:B6").PasteSpecial(-4122)
$q3.Range("D2:G6").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 3. Leave the selection where it was originally (last sheet, "2021-Q2").
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q2").Activate()
$wb.Worksheets.Item("2021-Q2").Range("A1").Select()

Write-Output "2022-Q3 sheet + summary row added"
